$wb = $excel.ActiveWorkbook

# --- Sheet2: replace the "Days" column (E2:E11) with combined Abc/Xyz labels ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Activate()

$newDays = @(
    "Abc1 Xyz1",
    "Abc2 Xyz2",
    "Abc3 Xyz3",
    "Abc4 Xyz4",
    "Abc5 Xyz5",
    "Abc6 Xyz6",
    "Abc7 Xyz7",
    "Abc8 Xyz8",
    "Abc9 Xyz9",
    "Abc10 Xyz10"
)

for ($i = 0; $i -lt $newDays.Length; $i++) {
    $row = $i + 2
    $ws2.Range("E$row").Value = $newDays[$i]
}

# Update the selection on Sheet2 to the refreshed Days column
$ws2.Range("E2:E11").Select()

# --- Sheet1: widen the selection to include column B ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("A2:B11").Select()

# --- Sheet3: move the selection down to A12:B21 ---
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.Range("A12:B21").Select()

# Re-activate Sheet2 as the last active tab (matches tabSelected="1" on Sheet2)
$ws2.Activate()
